# The commit swaps the embedded picture "name" attributes between the two
# logo images that live in the document's headers/footers:
#   - the Pearson/Edexcel PNG logo (currently "image2.png") becomes "image1.png"
#   - the BTEC JPG logo (currently "image1.jpg") becomes "image2.jpg"
# (the "descr" attribute / alt text is left untouched in both cases).
#
# These values live in <wp:docPr .../> and <pic:cNvPr .../> attributes inside
# the drawing XML, which is not part of the visible text Range.Find searches,
# so we rewrite the inline picture's own WordOpenXML fragment directly.

$d = $word.ActiveDocument

$renames = @{
    "image1.jpg" = "image2.jpg"
    "image2.png" = "image1.png"
}

function Update-LogoNames($rangeOwner) {
    if (-not $rangeOwner.Exists) { return }
    foreach ($shp in $rangeOwner.Range.InlineShapes) {
        $r = $shp.Range
        $xml = $r.WordOpenXML
        foreach ($old in $renames.Keys) {
            $new = $renames[$old]
            $needle = 'name="' + $old + '"'
            $replacement = 'name="' + $new + '"'
            if ($xml.Contains($needle)) {
                $xml = $xml.Replace($needle, $replacement)
            }
        }
        $r.WordOpenXML = $xml
    }
}

foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        Update-LogoNames $hf
    }
    foreach ($hf in $sec.Footers) {
        Update-LogoNames $hf
    }
}
